# Automated data refresh: updated extraction timestamps (DATA_EXTRACCIO)
# and a handful of revised observation readings (humidity, pressure,
# radiation, precipitation, temperature) for 2026-02-06.
#
# Note: values such as "87%" are written with a leading apostrophe so
# Excel stores them as literal text (matching the source column, which is
# plain text, not a numeric percentage) instead of auto-converting them to
# a numeric percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 22:47:39"
$ws.Range("E3").Value = "2026-02-06 22:47:42"
$ws.Range("E4").Value = "2026-02-06 22:47:44"
$ws.Range("K4").Value = "11.3 MJ/m2"
$ws.Range("O4").Value = "13.4 °C"
$ws.Range("E5").Value = "2026-02-06 22:47:46"
$ws.Range("E6").Value = "2026-02-06 22:47:49"
$ws.Range("E7").Value = "2026-02-06 22:47:51"
$ws.Range("E8").Value = "2026-02-06 22:47:54"
$ws.Range("O8").Value = "9.5 °C"
$ws.Range("E9").Value = "2026-02-06 22:47:56"
$ws.Range("H9").Value = "'87%"
$ws.Range("E10").Value = "2026-02-06 22:47:58"
$ws.Range("E11").Value = "2026-02-06 22:48:01"
$ws.Range("J11").Value = "999.6 hPa"
$ws.Range("O11").Value = "4.9 °C"
$ws.Range("E12").Value = "2026-02-06 22:48:03"
$ws.Range("O12").Value = "13.1 °C"
$ws.Range("E13").Value = "2026-02-06 22:48:06"
$ws.Range("O13").Value = "9.9 °C"
$ws.Range("E14").Value = "2026-02-06 22:48:08"
$ws.Range("H14").Value = "'76%"
$ws.Range("K14").Value = "7.5 MJ/m2"
$ws.Range("E15").Value = "2026-02-06 22:48:10"
$ws.Range("H15").Value = "'74%"
$ws.Range("O15").Value = "10.2 °C"
$ws.Range("E16").Value = "2026-02-06 22:48:13"
$ws.Range("O16").Value = "5.9 °C"
$ws.Range("E17").Value = "2026-02-06 22:48:15"
$ws.Range("J17").Value = "999.6 hPa"
$ws.Range("E18").Value = "2026-02-06 22:48:18"
$ws.Range("I18").Value = "0.8 mm"
$ws.Range("O18").Value = "-4.8 °C"
$ws.Range("E19").Value = "2026-02-06 22:48:20"
$ws.Range("H19").Value = "'79%"
$ws.Range("J19").Value = "1000.5 hPa"
$ws.Range("E20").Value = "2026-02-06 22:48:23"
$ws.Range("E21").Value = "2026-02-06 22:48:25"
$ws.Range("J21").Value = "998.5 hPa"
$ws.Range("O21").Value = "8.6 °C"
$ws.Range("E22").Value = "2026-02-06 22:48:28"
$ws.Range("O22").Value = "10.0 °C"
$ws.Range("E23").Value = "2026-02-06 22:48:30"
$ws.Range("E24").Value = "2026-02-06 22:48:32"
$ws.Range("O24").Value = "12.7 °C"
$ws.Range("E25").Value = "2026-02-06 22:48:35"
$ws.Range("J25").Value = "999.2 hPa"
$ws.Range("O25").Value = "4.3 °C"
$ws.Range("E26").Value = "2026-02-06 22:48:37"
$ws.Range("O26").Value = "-1.1 °C"
$ws.Range("E27").Value = "2026-02-06 22:48:39"
$ws.Range("H27").Value = "'79%"
$ws.Range("J27").Value = "998.3 hPa"
$ws.Range("E28").Value = "2026-02-06 22:48:42"
$ws.Range("H28").Value = "'84%"
$ws.Range("J28").Value = "1000.5 hPa"
$ws.Range("O28").Value = "4.9 °C"
$ws.Range("E29").Value = "2026-02-06 22:48:44"
$ws.Range("O29").Value = "12.1 °C"
$ws.Range("E30").Value = "2026-02-06 22:48:47"
$ws.Range("E31").Value = "2026-02-06 22:48:49"
$ws.Range("I31").Value = "1.9 mm"
$ws.Range("J31").Value = "1000.0 hPa"
$ws.Range("O31").Value = "7.1 °C"
$ws.Range("E32").Value = "2026-02-06 22:48:52"
$ws.Range("J32").Value = "999.5 hPa"
$ws.Range("E33").Value = "2026-02-06 22:48:54"
$ws.Range("O33").Value = "10.1 °C"
$ws.Range("E34").Value = "2026-02-06 22:48:56"
$ws.Range("E35").Value = "2026-02-06 22:48:59"
$ws.Range("N35").Value = "-4.0 °C 22:27 TU"
$ws.Range("O35").Value = "-2.3 °C"
$ws.Range("E36").Value = "2026-02-06 22:49:01"
$ws.Range("H36").Value = "'69%"
$ws.Range("J36").Value = "1000.6 hPa"
$ws.Range("N36").Value = "6.8 °C 22:27 TU"
$ws.Range("O36").Value = "12.3 °C"
